$d = $word.ActiveDocument

$d.Content.Find.Execute("35÷8=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷3=13, 0", 2) | Out-Null
$d.Content.Find.Execute("53÷3=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "68÷6=11, 2", 2) | Out-Null
$d.Content.Find.Execute("81÷8=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷3=27, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷9=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "36÷4=9, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷6=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "71÷9=7, 8", 2) | Out-Null
$d.Content.Find.Execute("19÷4=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "48÷2=24, 0", 2) | Out-Null
$d.Content.Find.Execute("35÷9=3, 8", $true, $false, $false, $false, $false, $true, 1, $false, "20÷7=2, 6", 2) | Out-Null
$d.Content.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "45÷3=15, 0", 2) | Out-Null
$d.Content.Find.Execute("93÷2=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
$d.Content.Find.Execute("23÷7=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=4, 3", 2) | Out-Null
$d.Content.Find.Execute("75÷6=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "95÷5=19, 0", 2) | Out-Null
$d.Content.Find.Execute("95÷6=15, 5", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷4=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "42÷8=5, 2", 2) | Out-Null
$d.Content.Find.Execute("57÷8=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=8, 4", 2) | Out-Null
$d.Content.Find.Execute("90÷4=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=43, 0", 2) | Out-Null
$d.Content.Find.Execute("91÷5=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=6, 0", 2) | Out-Null
$d.Content.Find.Execute("64÷4=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷5=16, 2", 2) | Out-Null
$d.Content.Find.Execute("96÷5=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=15, 0", 2) | Out-Null
$d.Content.Find.Execute("30÷7=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷9=5, 5", 2) | Out-Null
$d.Content.Find.Execute("54÷7=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "35÷6=5, 5", 2) | Out-Null
$d.Content.Find.Execute("87÷3=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "64÷9=7, 1", 2) | Out-Null
$d.Content.Find.Execute("65÷5=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=10, 6", 2) | Out-Null
$d.Content.Find.Execute("58÷3=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "13÷6=2, 1", 2) | Out-Null
$d.Content.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "64÷8=8, 0", 2) | Out-Null
$d.Content.Find.Execute("97÷6=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "64÷8=8, 0", 2) | Out-Null
